$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.426.36'
$ws.Range("E2").Value = '  -2.41%  '

$ws.Range("D3").Value = '2.618.33'
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.33'
$ws.Range("E5").Value = '  -0.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.48'
$ws.Range("E6").Value = '  -4.36%  '

$ws.Range("E7").Value = '  -1.46%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("E9").Value = '  -0.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.50'
$ws.Range("E10").Value = '  +0.07%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0845'
$ws.Range("E11").Value = '  +0.35%  '

$ws.Range("B12").Value = 'OKB'
$ws.Range("C12").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.31'
$ws.Range("E12").Value = '  -1.33%  '

$ws.Range("E13").Value = '  +0.98%  '

$ws.Range("D14").Value = '3.016.11'
$ws.Range("E14").Value = '  -0.29%  '

$ws.Range("E15").Value = '  +0.37%  '

$ws.Range("D16").Value = '2.622.05'
$ws.Range("E16").Value = '  +1.37%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.927'
$ws.Range("E17").Value = '  +1.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.02'
$ws.Range("E18").Value = '  -2.05%  '

$ws.Range("D19").Value = '46.596.60'
$ws.Range("E19").Value = '  -2.34%  '

$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("E21").Value = '  -7.70%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.76'
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.87'
$ws.Range("E23").Value = '  +1.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '275.57'
$ws.Range("E24").Value = '  +6.19%  '

$ws.Range("E25").Value = '  +0.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.20'
$ws.Range("E26").Value = '  +2.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.55'
$ws.Range("E27").Value = '  +12.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.16%  '

$ws.Range("E29").Value = '  -1.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.65'
$ws.Range("E30").Value = '  +0.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '38.72'
$ws.Range("E31").Value = '  -8.94%  '

$ws.Range("E32").Value = '  -2.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.46'
$ws.Range("E33").Value = '  +6.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.63'
$ws.Range("E34").Value = '  -6.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.24'
$ws.Range("E35").Value = '  -0.77%  '

$ws.Range("E36").Value = '  -5.81%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0841'
$ws.Range("E37").Value = '  -1.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.36'
$ws.Range("E38").Value = '  +0.98%  '

$ws.Range("E39").Value = '  -0.54%  '

$ws.Range("E40").Value = '  +1.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.47'
$ws.Range("E41").Value = '  +29.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.98'
$ws.Range("E42").Value = '  -4.91%  '

$ws.Range("E43").Value = '  -0.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.64'
$ws.Range("E44").Value = '  -0.84%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.06'
$ws.Range("E45").Value = '  -6.62%  '

$ws.Range("D46").Value = '2.133.26'
$ws.Range("E46").Value = '  +5.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("E47").Value = '  +0.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '94.39'
$ws.Range("E48").Value = '  -0.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.54'
$ws.Range("E49").Value = '  +6.00%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.05'
$ws.Range("E50").Value = '  +0.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.79'
$ws.Range("E51").Value = '  -7.76%  '
